
# Update profit files after running on 2025-10-17
# Appends a new data row (row 61) to the sheet with the date "10/17/2025"
# and profit value 9501.02, extending the data range from A1:B60 to A1:B61.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 61

# The Date column stores dates as literal text (e.g. "10/16/2025") rather
# than true date serial values, matching the existing rows above it.
# Force the cell to Text format before assigning the string so Excel does
# not auto-convert the "mm/dd/yyyy"-looking string into a date value, then
# clear the explicit formatting again so the cell keeps the sheet's default
# (unstyled) appearance, just like the other plain date cells.
$dateCell = $ws.Cells.Item($newRow, 1)
$dateCell.NumberFormat = "@"
$dateCell.Value = "10/17/2025"
$dateCell.ClearFormats()

# Profit column holds a plain numeric value.
$ws.Cells.Item($newRow, 2).Value = 9501.02
